# DRILL-8006: Leading and Trailing Whitespace Causes Query Failures in Excel Files
#
# Adds a new worksheet ("spaceInColHeader") after the existing "comps" sheet
# that exercises a column header with a trailing space ("col2 "), used to
# reproduce/verify the whitespace-handling fix described in the commit.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "comps" (it becomes the last / newly
# active sheet, matching the workbook's activeTab pointing at it afterwards).
$comps = $wb.Worksheets.Item("comps")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $comps)
$newSheet.Name = "spaceInColHeader"

# Header row: "col1" is a plain header, "col2 " intentionally keeps a
# trailing space to reproduce the whitespace bug.
$newSheet.Range("A1").Value = "col1"
$newSheet.Range("B1").Value = "col2 "

# A couple of simple data rows under the headers.
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("A3").Value = 3
$newSheet.Range("B3").Value = 4

# Make the new sheet the active tab/selection (B1), mirroring the saved
# workbook state where this is the last-edited sheet.
$newSheet.Activate() | Out-Null
$newSheet.Range("B1").Select() | Out-Null
